$d = $word.ActiveDocument

$d.Content.Find.Execute("67÷7=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "46÷6=7, 4", 2) | Out-Null
$d.Content.Find.Execute("45÷3=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "19÷4=4, 3", 2) | Out-Null
$d.Content.Find.Execute("51÷2=25, 1", $true, $false, $false, $false, $false, $true, 1, $false, "92÷6=15, 2", 2) | Out-Null
$d.Content.Find.Execute("96÷8=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "74÷8=9, 2", 2) | Out-Null
$d.Content.Find.Execute("17÷4=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "56÷7=8, 0", 2) | Out-Null
$d.Content.Find.Execute("70÷9=7, 7", $true, $false, $false, $false, $false, $true, 1, $false, "53÷2=26, 1", 2) | Out-Null
$d.Content.Find.Execute("35÷6=5, 5", $true, $false, $false, $false, $false, $true, 1, $false, "30÷8=3, 6", 2) | Out-Null
$d.Content.Find.Execute("83÷3=27, 2", $true, $false, $false, $false, $false, $true, 1, $false, "96÷3=32, 0", 2) | Out-Null
$d.Content.Find.Execute("56÷5=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "91÷4=22, 3", 2) | Out-Null
$d.Content.Find.Execute("32÷7=4, 4", $true, $false, $false, $false, $false, $true, 1, $false, "88÷7=12, 4", 2) | Out-Null
$d.Content.Find.Execute("32÷8=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "16÷3=5, 1", 2) | Out-Null
$d.Content.Find.Execute("47÷8=5, 7", $true, $false, $false, $false, $false, $true, 1, $false, "42÷8=5, 2", 2) | Out-Null
$d.Content.Find.Execute("62÷5=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "98÷4=24, 2", 2) | Out-Null
$d.Content.Find.Execute("96÷2=48, 0", $true, $false, $false, $false, $false, $true, 1, $false, "59÷3=19, 2", 2) | Out-Null
$d.Content.Find.Execute("90÷2=45, 0", $true, $false, $false, $false, $false, $true, 1, $false, "29÷8=3, 5", 2) | Out-Null
$d.Content.Find.Execute("35÷2=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "66÷6=11, 0", 2) | Out-Null
$d.Content.Find.Execute("83÷5=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "49÷3=16, 1", 2) | Out-Null
$d.Content.Find.Execute("75÷5=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "32÷5=6, 2", 2) | Out-Null
$d.Content.Find.Execute("64÷6=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "31÷2=15, 1", 2) | Out-Null
$d.Content.Find.Execute("31÷3=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "60÷9=6, 6", 2) | Out-Null
$d.Content.Find.Execute("57÷3=19, 0", $true, $false, $false, $false, $false, $true, 1, $false, "29÷3=9, 2", 2) | Out-Null
$d.Content.Find.Execute("36÷7=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "45÷7=6, 3", 2) | Out-Null
$d.Content.Find.Execute("49÷8=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "54÷9=6, 0", 2) | Out-Null
$d.Content.Find.Execute("18÷5=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "94÷8=11, 6", 2) | Out-Null
$d.Content.Find.Execute("60÷4=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "77÷9=8, 5", 2) | Out-Null
